$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- style bookkeeping: Text + Date-time number formats (with a bordered
# edge each) get created once, mirroring formatting experiments the author
# ran on this sheet; the scratch cells are cleared afterward so only the
# style table entries remain. ---
$ws.Range("Z1").Borders.Item(7).LineStyle = 1
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z2").Borders.Item(8).LineStyle = 1
$ws.Range("Z2").NumberFormat = "m/d/yy h:mm"
$ws.Range("Z1:Z2").ClearFormats()
$ws.Range("Z1:Z2").ClearContents()

# --- column widths ---
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 12
$ws.Columns.Item(8).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 15
$ws.Columns.Item(10).ColumnWidth = 12
$ws.Columns.Item(11).ColumnWidth = 12
$ws.Columns.Item(12).ColumnWidth = 12

# --- YCbCr 135-degree GLCM feature data ---
$ws.Range("A1").Value = 0.13735896374203449
$ws.Range("B1").Value = 0.71632427367301299
$ws.Range("C1").Value = 0.55219022973870036
$ws.Range("D1").Value = 0.935598965647223
$ws.Range("E1").Value = 0.0087603365023693723
$ws.Range("F1").Value = 0.78281185019536137
$ws.Range("G1").Value = 0.94072490532360598
$ws.Range("H1").Value = 0.99561988151785152
$ws.Range("I1").Value = 0.0028705954695445197
$ws.Range("J1").Value = 0.68832114594817606
$ws.Range("K1").Value = 0.98512121346181902
$ws.Range("L1").Value = 0.99856470226522775
$ws.Range("A2").Value = 0.15309602437335543
$ws.Range("B2").Value = 0.78465446497925861
$ws.Range("C2").Value = 0.47460500301462288
$ws.Range("D2").Value = 0.92643505210429755
$ws.Range("E2").Value = 0.028781923262434797
$ws.Range("F2").Value = 0.80190703111331285
$ws.Range("G2").Value = 0.78897385733639125
$ws.Range("H2").Value = 0.9856091658158086
$ws.Range("I2").Value = 0.00046820655957754829
$ws.Range("J2").Value = 0.16717507726009193
$ws.Range("K2").Value = 0.99888776491555198
$ws.Range("L2").Value = 0.99976589672021121
$ws.Range("A3").Value = 0.23038861128720584
$ws.Range("B3").Value = 0.66079378668775923
$ws.Range("C3").Value = 0.38160838928927582
$ws.Range("D3").Value = 0.89025959657689202
$ws.Range("E3").Value = 0.077844955933743951
$ws.Range("F3").Value = 0.82150457728695925
$ws.Range("G3").Value = 0.49625887081502845
$ws.Range("H3").Value = 0.96107752203312835
$ws.Range("I3").Value = 0.00015862346977564033
$ws.Range("J3").Value = 0.69807130911463189
$ws.Range("K3").Value = 0.99928543789068036
$ws.Range("L3").Value = 0.99992068826511216
$ws.Range("A4").Value = 0.13595643975591323
$ws.Range("B4").Value = 0.8280506040927198
$ws.Range("C4").Value = 0.34089024450013533
$ws.Range("D4").Value = 0.93215082228894386
$ws.Range("E4").Value = 0.02777273630679325
$ws.Range("F4").Value = 0.86927619585787408
$ws.Range("G4").Value = 0.76506674371925798
$ws.Range("H4").Value = 0.98611363184660328
$ws.Range("I4").Value = 0.061628073840773012
$ws.Range("J4").Value = 0.79679556742922075
$ws.Range("K4").Value = 0.64221301870417624
$ws.Range("L4").Value = 0.96918596307961347
